# Trade #4 closed at 2026-02-16 22:56:44 - base_strategy DOWN +0.000%
# Appends a new trade row (row 5) to both the "All Trades" and
# "base_strategy" worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $row = 5

    # Trade #
    $ws.Cells.Item($row, 1).Value = 4

    # Date - force text so "2026-02-16" isn't auto-converted to a date serial
    $cDate = $ws.Cells.Item($row, 2)
    $cDate.NumberFormat = "@"
    $cDate.Value = "2026-02-16"
    $cDate.Style = "Normal"

    # Time
    $ws.Cells.Item($row, 3).Value = "22:56:44"

    # Strategy
    $ws.Cells.Item($row, 4).Value = "base_strategy"

    # Side
    $ws.Cells.Item($row, 5).Value = "DOWN"

    # Entry Price
    $ws.Cells.Item($row, 6).Value = 0.5

    # Exit Price - stays blank (present empty cell, like the rest of column G)
    $cExit = $ws.Cells.Item($row, 7)
    $cExit.Style = "Normal"

    # Status
    $ws.Cells.Item($row, 8).Value = "OPEN"

    # P&L %
    $ws.Cells.Item($row, 9).Value = 0

    # P&L $
    $ws.Cells.Item($row, 10).Value = 0

    # Capital After
    $ws.Cells.Item($row, 11).Value = 100

    # Entry Slippage (bps)
    $ws.Cells.Item($row, 12).Value = 0

    # Exit Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0

    # Confidence
    $ws.Cells.Item($row, 14).Value = 0.6

    # Entry Reason
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"

    # Exit Reason - stays blank (present empty cell, like the rest of column P)
    $cExitReason = $ws.Cells.Item($row, 16)
    $cExitReason.Style = "Normal"

    # Duration (min)
    $ws.Cells.Item($row, 17).Value = 0
}
